$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45041
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11500
$ws.Range("Q2").Value = '$/bandeja 18 kilos granel'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 639
$ws.Range("T2").Value = 18

$ws.Range("D3").Value = 45050
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 140
$ws.Range("N3").Value = 11000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 11429
$ws.Range("Q3").Value = '$/caja 18 kilos empedrada'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 635
$ws.Range("T3").Value = 18

$ws.Range("D4").Value = 45029
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 9000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 9500
$ws.Range("Q4").Value = '$/bandeja 18 kilos granel'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 528
$ws.Range("T4").Value = 18

$ws.Range("D5").Value = 45013
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 9000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 9500
$ws.Range("Q5").Value = '$/bandeja 18 kilos granel'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 528
$ws.Range("T5").Value = 18

$ws.Range("D6").Value = 45079
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 270
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 11444
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = 'Región de O''Higgins'
$ws.Range("S6").Value = 636
$ws.Range("T6").Value = 18

$ws.Range("D7").Value = 44272
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 9000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 9500
$ws.Range("Q7").Value = '$/caja 15 kilos granel'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 633
$ws.Range("T7").Value = 15

$ws.Range("D8").Value = 44272
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("Q8").Value = '$/caja 15 kilos granel'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 533
$ws.Range("T8").Value = 15

$ws.Range("D9").Value = 45126
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("Q9").Value = '$/bandeja 18 kilos granel'
$ws.Range("R9").Value = 'Región de O''Higgins'
$ws.Range("S9").Value = 806
$ws.Range("T9").Value = 18

$ws.Range("D10").Value = 45034
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 220
$ws.Range("N10").Value = 8500
$ws.Range("O10").Value = 9000
$ws.Range("P10").Value = 8727
$ws.Range("Q10").Value = '$/caja 18 kilos granel'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 485
$ws.Range("T10").Value = 18

$ws.Range("D11").Value = 44307
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 50
$ws.Range("N11").Value = 10000
$ws.Range("O11").Value = 10000
$ws.Range("P11").Value = 10000
$ws.Range("Q11").Value = '$/bandeja 18 kilos granel'
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 556
$ws.Range("T11").Value = 18

$ws.Range("D12").Value = 44307
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 8000
$ws.Range("O12").Value = 8000
$ws.Range("P12").Value = 8000
$ws.Range("Q12").Value = '$/bandeja 18 kilos granel'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 444
$ws.Range("T12").Value = 18

$ws.Range("D13").Value = 44698
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 10000
$ws.Range("O13").Value = 10000
$ws.Range("P13").Value = 10000
$ws.Range("Q13").Value = '$/caja 18 kilos granel'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 556
$ws.Range("T13").Value = 18

$ws.Range("D14").Value = 45014
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 9000
$ws.Range("O14").Value = 10000
$ws.Range("P14").Value = 9500
$ws.Range("Q14").Value = '$/bandeja 18 kilos granel'
$ws.Range("R14").Value = 'Región de O''Higgins'
$ws.Range("S14").Value = 528
$ws.Range("T14").Value = 18

$ws.Range("D15").Value = 45128
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 12000
$ws.Range("O15").Value = 12000
$ws.Range("P15").Value = 12000
$ws.Range("Q15").Value = '$/bandeja 18 kilos granel'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 667
$ws.Range("T15").Value = 18

$ws.Range("D16").Value = 44358
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 11000
$ws.Range("O16").Value = 12000
$ws.Range("P16").Value = 11500
$ws.Range("Q16").Value = '$/caja 18 kilos granel'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 639
$ws.Range("T16").Value = 18

$ws.Range("D17").Value = 45107
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 11000
$ws.Range("O17").Value = 11000
$ws.Range("P17").Value = 11000
$ws.Range("Q17").Value = '$/caja 18 kilos empedrada'
$ws.Range("R17").Value = 'Región del Maule'
$ws.Range("S17").Value = 611
$ws.Range("T17").Value = 18

$ws.Range("D18").Value = 44425
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 12000
$ws.Range("O18").Value = 13000
$ws.Range("P18").Value = 12500
$ws.Range("Q18").Value = '$/bandeja 18 kilos granel'
$ws.Range("R18").Value = 'Región de O''Higgins'
$ws.Range("S18").Value = 694
$ws.Range("T18").Value = 18

$ws.Range("D19").Value = 45076
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 150
$ws.Range("N19").Value = 10000
$ws.Range("O19").Value = 11000
$ws.Range("P19").Value = 10467
$ws.Range("Q19").Value = '$/caja 18 kilos granel'
$ws.Range("R19").Value = 'Provincia de Curicó'
$ws.Range("S19").Value = 582
$ws.Range("T19").Value = 18

$ws.Range("D20").Value = 45154
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 13000
$ws.Range("O20").Value = 14000
$ws.Range("P20").Value = 13500
$ws.Range("Q20").Value = '$/bandeja 18 kilos granel'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 750
$ws.Range("T20").Value = 18

$ws.Range("D21").Value = 44316
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 9000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 9500
$ws.Range("Q21").Value = '$/caja 18 kilos granel'
$ws.Range("R21").Value = 'Región de O''Higgins'
$ws.Range("S21").Value = 528
$ws.Range("T21").Value = 18

$ws.Range("D22").Value = 44363
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 9000
$ws.Range("O22").Value = 10000
$ws.Range("P22").Value = 9500
$ws.Range("Q22").Value = '$/caja 15 kilos empedrada'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 633
$ws.Range("T22").Value = 15

$ws.Range("D23").Value = 44999
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 12000
$ws.Range("O23").Value = 12000
$ws.Range("P23").Value = 12000
$ws.Range("Q23").Value = '$/bandeja 18 kilos granel'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 667
$ws.Range("T23").Value = 18

$ws.Range("D24").Value = 44999
$ws.Range("L24").Value = 'Segunda'
$ws.Range("M24").Value = 100
$ws.Range("N24").Value = 10000
$ws.Range("O24").Value = 10000
$ws.Range("P24").Value = 10000
$ws.Range("Q24").Value = '$/bandeja 18 kilos granel'
$ws.Range("R24").Value = 'Región de O''Higgins'
$ws.Range("S24").Value = 556
$ws.Range("T24").Value = 18

$ws.Range("D25").Value = 44776
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 10000
$ws.Range("O25").Value = 10000
$ws.Range("P25").Value = 10000
$ws.Range("Q25").Value = '$/bandeja 18 kilos granel'
$ws.Range("R25").Value = 'Región de O''Higgins'
$ws.Range("S25").Value = 556
$ws.Range("T25").Value = 18

$ws.Range("D26").Value = 44776
$ws.Range("L26").Value = 'Segunda'
$ws.Range("M26").Value = 50
$ws.Range("N26").Value = 8000
$ws.Range("O26").Value = 8000
$ws.Range("P26").Value = 8000
$ws.Range("Q26").Value = '$/bandeja 18 kilos granel'
$ws.Range("R26").Value = 'Región de O''Higgins'
$ws.Range("S26").Value = 444
$ws.Range("T26").Value = 18

$ws.Range("D27").Value = 44299
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 10000
$ws.Range("O27").Value = 11000
$ws.Range("P27").Value = 10500
$ws.Range("Q27").Value = '$/caja 18 kilos granel'
$ws.Range("R27").Value = 'Región del Maule'
$ws.Range("S27").Value = 583
$ws.Range("T27").Value = 18

$ws.Range("D28").Value = 44299
$ws.Range("L28").Value = 'Segunda'
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 9000
$ws.Range("O28").Value = 9000
$ws.Range("P28").Value = 9000
$ws.Range("Q28").Value = '$/caja 18 kilos granel'
$ws.Range("R28").Value = 'Región del Maule'
$ws.Range("S28").Value = 500
$ws.Range("T28").Value = 18

$ws.Range("D29").Value = 45140
$ws.Range("L29").Value = 'Primera'
$ws.Range("M29").Value = 50
$ws.Range("N29").Value = 12000
$ws.Range("O29").Value = 12000
$ws.Range("P29").Value = 12000
$ws.Range("Q29").Value = '$/bandeja 18 kilos granel'
$ws.Range("R29").Value = 'Región de O''Higgins'
$ws.Range("S29").Value = 667
$ws.Range("T29").Value = 18

$ws.Range("D30").Value = 45092
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 110
$ws.Range("N30").Value = 10000
$ws.Range("O30").Value = 11000
$ws.Range("P30").Value = 10455
$ws.Range("Q30").Value = '$/bandeja 18 kilos granel'
$ws.Range("R30").Value = 'Provincia de Curicó'
$ws.Range("S30").Value = 581
$ws.Range("T30").Value = 18

$ws.Range("D31").Value = 45027
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 100
$ws.Range("N31").Value = 9000
$ws.Range("O31").Value = 10000
$ws.Range("P31").Value = 9500
$ws.Range("Q31").Value = '$/bandeja 18 kilos granel'
$ws.Range("R31").Value = 'Región de O''Higgins'
$ws.Range("S31").Value = 528
$ws.Range("T31").Value = 18

$ws.Range("D32").Value = 45037
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 250
$ws.Range("N32").Value = 9000
$ws.Range("O32").Value = 9500
$ws.Range("P32").Value = 9200
$ws.Range("Q32").Value = '$/caja 18 kilos granel'
$ws.Range("R32").Value = 'Provincia de Curicó'
$ws.Range("S32").Value = 511
$ws.Range("T32").Value = 18

